$d = $word.ActiveDocument

# Merge the run boundaries by replacing the text spanning the first three
# runs (ending right before "implementation.") with the identical text,
# so Word collapses it into a single run. Likewise for the run spanning
# from "implementation." through "...architectures,".

$d.Content.Find.Execute(
    "database management systems(DBMS); as well as, data analysis, database design, data modeling, database management and database",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "database management systems(DBMS); as well as, data analysis, database design, data modeling, database management and database",
    2
) | Out-Null

$d.Content.Find.Execute(
    "implementation. There is a specific emphasis on data analytics and learning to query data with Structured Query Language (SQL), query performance, data normalization; and database migration. This course provides hands-on experience in database design and implementation through assignments, lab exercises and course projects. This course also introduces advanced database concepts such as transaction management and concurrency control, distributed databases, multi-tier client/server architectures,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "implementation. There is a specific emphasis on data analytics and learning to query data with Structured Query Language (SQL), query performance, data normalization; and database migration. This course provides hands-on experience in database design and implementation through assignments, lab exercises and course projects. This course also introduces advanced database concepts such as transaction management and concurrency control, distributed databases, multi-tier client/server architectures,",
    2
) | Out-Null
